$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "RFCU4114542"
$ws.Range("B5").Value = "EVER LEGACY"
$ws.Range("C5").Value = "034E"
$ws.Range("D5").Value = "9072904799-01"
$ws.Range("E5").Font.Bold = $false
$ws.Range("F5").Value = "COSU6203869480"

$ws.Range("A6").Value = "HESU4027089"
$ws.Range("B6").Value = "NAVARINO"
$ws.Range("C6").Value = "0848E"
$ws.Range("D6").Value = "9075904134-01"
$ws.Range("E6").Font.Bold = $false
$ws.Range("F6").Value = "COSU6203956310"
